$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price cells keep their text formatting so values such as
# "29.360.13" or "1.000" are not reinterpreted as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value2 = '29.360.13'
$ws.Range("E2").Value2 = '  +0.12%  '
$ws.Range("D3").Value2 = '1.874.60'
$ws.Range("E3").Value2 = '  +0.22%  '
$ws.Range("D4").Value2 = '1.000'
$ws.Range("E4").Value2 = '  -0.14%  '
$ws.Range("D5").Value2 = '0.7126'
$ws.Range("E5").Value2 = '  -1.09%  '
$ws.Range("D6").Value2 = '241.60'
$ws.Range("E6").Value2 = '  +0.20%  '
$ws.Range("E7").Value2 = '  -0.09%  '
$ws.Range("D8").Value2 = '0.3109'
$ws.Range("E8").Value2 = '  +0.57%  '
$ws.Range("D9").Value2 = '0.07690'
$ws.Range("E9").Value2 = '  -2.33%  '
$ws.Range("D10").Value2 = '25.17'
$ws.Range("E10").Value2 = '  -0.55%  '
$ws.Range("D11").Value2 = '0.08372'
$ws.Range("E11").Value2 = '  +1.22%  '
$ws.Range("D12").Value2 = '1.882.67'
$ws.Range("E12").Value2 = '  +1.39%  '
$ws.Range("D13").Value2 = '5.239'
$ws.Range("E13").Value2 = '  -0.12%  '
$ws.Range("D14").Value2 = '0.7140'
$ws.Range("E14").Value2 = '  -1.15%  '
$ws.Range("D15").Value2 = '91.58'
$ws.Range("E15").Value2 = '  +0.92%  '
$ws.Range("B16").Value2 = 'ShibaInu'
$ws.Range("C16").Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value2 = '0.000008349'
$ws.Range("E16").Value2 = '  +6.70%  '
$ws.Range("B17").Value2 = 'WrappedBTC'
$ws.Range("C17").Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value2 = '29.375.19'
$ws.Range("E17").Value2 = '  +0.16%  '
$ws.Range("D18").Value2 = '5.963'
$ws.Range("E18").Value2 = '  +1.94%  '
$ws.Range("D19").Value2 = '243.41'
$ws.Range("E19").Value2 = '  -0.12%  '
$ws.Range("D20").Value2 = '2.133.53'
$ws.Range("E20").Value2 = '  +1.26%  '
$ws.Range("D21").Value2 = '13.19'
$ws.Range("E21").Value2 = '  -0.30%  '
$ws.Range("E22").Value2 = '  -0.15%  '
$ws.Range("D23").Value2 = '7.893'
$ws.Range("E23").Value2 = '  -1.17%  '
$ws.Range("E24").Value2 = '  -0.17%  '
$ws.Range("D25").Value2 = '0.1614'
$ws.Range("E25").Value2 = '  +0.19%  '
$ws.Range("D26").Value2 = '163.97'
$ws.Range("E26").Value2 = '  +0.90%  '
$ws.Range("D27").Value2 = '9.001'
$ws.Range("E27").Value2 = '  +0.36%  '
$ws.Range("D28").Value2 = '18.56'
$ws.Range("E28").Value2 = '  +1.67%  '
$ws.Range("D29").Value2 = '1.506'
$ws.Range("E29").Value2 = '  +0.71%  '
$ws.Range("D30").Value2 = '4.403'
$ws.Range("E30").Value2 = '  +0.64%  '
$ws.Range("B31").Value2 = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value2 = '4.327'
$ws.Range("E31").Value2 = '  +5.28%  '
$ws.Range("B32").Value2 = 'Toncoin'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value2 = '1.290'
$ws.Range("E32").Value2 = '  -4.49%  '
$ws.Range("D33").Value2 = '0.05231'
$ws.Range("E33").Value2 = '  +0.39%  '
$ws.Range("D34").Value2 = '1.920'
$ws.Range("E34").Value2 = '  -0.95%  '
$ws.Range("D35").Value2 = '0.7622'
$ws.Range("E35").Value2 = '  +4.71%  '
$ws.Range("D36").Value2 = '1.172'
$ws.Range("E36").Value2 = '  -1.21%  '
$ws.Range("E37").Value2 = '  -0.16%  '
$ws.Range("D38").Value2 = '0.01862'
$ws.Range("E38").Value2 = '  +0.05%  '
$ws.Range("D39").Value2 = '2.723'
$ws.Range("E39").Value2 = '  +0.83%  '
$ws.Range("D40").Value2 = '1.161.76'
$ws.Range("E40").Value2 = '  -1.33%  '
$ws.Range("D41").Value2 = '6.399'
$ws.Range("E41").Value2 = '  +4.54%  '
$ws.Range("D42").Value2 = '73.32'
$ws.Range("E42").Value2 = '  +1.13%  '
$ws.Range("D43").Value2 = '0.8884'
$ws.Range("E43").Value2 = '  -1.47%  '
$ws.Range("D44").Value2 = '104.56'
$ws.Range("E44").Value2 = '  +2.89%  '
$ws.Range("D45").Value2 = '1.000'
$ws.Range("E45").Value2 = '  -0.12%  '
$ws.Range("D46").Value2 = '2.029.96'
$ws.Range("E46").Value2 = '  +1.06%  '
$ws.Range("D47").Value2 = '0.5202'
$ws.Range("E47").Value2 = '  -1.63%  '
$ws.Range("D48").Value2 = '1.794'
$ws.Range("E48").Value2 = '  +0.47%  '
$ws.Range("B49").Value2 = 'BabyDogeCoin'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value2 = '0.00000000120'
$ws.Range("E49").Value2 = '  -0.67%  '
$ws.Range("B50").Value2 = 'EnergySwap'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value2 = '9.391'
$ws.Range("E50").Value2 = '  +1.28%  '
$ws.Range("D51").Value2 = '0.4300'
$ws.Range("E51").Value2 = '  +0.38%  '
